# The daily auto-push job appends a new timestamp row for "today" (2026/01/27)
# right after the last row that already carried that date, shifting the
# remaining rows (the pre-generated future calendar entries) down by one.
#
# Original layout around the edit point:
#   ...
#   row 721: 2026/01/27  火   1   201   <- last existing "today" row
#   row 722: 2026/12/29  火  13   201   <- first of the pre-generated future rows
#   ...
#   row 763: 2027/01/05  火   7   201   <- last row in the sheet
#
# New layout: a fresh "2026/01/27" sample (hour 5) is inserted at row 722,
# pushing everything that used to start at row 722 down to row 723, and the
# sheet grows by one row (new last row: 764).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 722 (shifts existing rows 722:763 down to 723:764).
$ws.Rows.Item(722).Insert()

# Column A holds dates stored as plain text (e.g. "2026/01/27"), not real
# Excel dates, in the rest of the sheet. Force text formatting before the
# assignment so the date-shaped string isn't auto-converted to a date
# serial number, then clear the formatting residue so the cell ends up
# styled exactly like its neighbours (no explicit style index).
$cellA = $ws.Cells.Item(722, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2026/01/27"
$cellA.ClearFormats()

$ws.Cells.Item(722, 2).Value = "火"
$ws.Cells.Item(722, 3).Value = 5
$ws.Cells.Item(722, 4).Value = 201
